$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8721416592597961
$ws.Range("B1").Value = 1.273996591567993
$ws.Range("C1").Value = 2.336717844009399
$ws.Range("D1").Value = 2.523666858673096
$ws.Range("E1").Value = 1.938113570213318
